$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lesson #16 (row 19): video material uploaded to YouTube - add the multi-part
# link text into the "YouTube link" column and wrap it like the other
# multi-line link cells in this sheet.
$newText = "Part #1: https://youtu.be/Adstn_hqvJg`nPart #2: https://youtu.be/HbMFNIIB-V4`nPart #3: https://youtu.be/CzmM2m1ouaI`nPart #4: https://youtu.be/ZWD2dJoyuGQ"
$ws.Range("F19").Value = $newText
$ws.Range("F19").WrapText = $true

# Lesson #17 (row 20): now has a date - copy the date formatting used by the
# cell above (E19) and fill in the actual date.
$ws.Range("E19").Copy() | Out-Null
$ws.Range("E20").PasteSpecial(-4122) | Out-Null
$ws.Range("E20").Value = 44158

# Row heights adjust to fit the newly-wrapped 4-line text / the restored
# single-line row below it.
$ws.Rows.Item(19).RowHeight = 55.2
$ws.Rows.Item(20).RowHeight = 14.9

# Move the sheet's active selection.
$ws.Range("F25").Select() | Out-Null
